$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.215.36"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "1.582.82"
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").Value = "'209.69"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("D6").Value = "'0.497"
$ws.Range("E6").Value = "  -2.99%  "

$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("E8").Value = "  -1.40%  "

$ws.Range("D9").Value = "'0.246"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").Value = "'19.51"
$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").Value = "1.805.07"
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.06"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.579.65"
$ws.Range("E14").Value = "  -1.47%  "

$ws.Range("D15").Value = "'0.516"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").Value = "'64.55"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "26.202.90"
$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").Value = "'7.27"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'207.99"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").Value = "'4.26"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("E23").Value = "  -3.45%  "

$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  -0.75%  "

$ws.Range("D25").Value = "'144.79"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").Value = "'7.02"
$ws.Range("E27").Value = "  -0.80%  "

$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").Value = "'15.21"
$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("D30").Value = "'0.0505"
$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("E31").Value = "  -1.41%  "

$ws.Range("E32").Value = "  -1.05%  "

$ws.Range("D33").Value = "'2.96"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("D34").Value = "1.280.22"
$ws.Range("E34").Value = "  -1.10%  "

$ws.Range("D35").Value = "'2.46"
$ws.Range("E35").Value = "  -0.46%  "

$ws.Range("D36").Value = "'0.609"
$ws.Range("E36").Value = "  +2.34%  "

$ws.Range("D37").Value = "'1.49"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").Value = "'1.14"
$ws.Range("E38").Value = "  -1.15%  "

$ws.Range("E39").Value = "  -1.51%  "

$ws.Range("D40").Value = "'0.815"
$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("E41").Value = "  +3.37%  "

$ws.Range("D42").Value = "'0.765"
$ws.Range("E42").Value = "  -2.08%  "

$ws.Range("E43").Value = "  -3.11%  "

$ws.Range("D44").Value = "'62.32"
$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").Value = "1.718.38"
$ws.Range("E45").Value = "  -1.19%  "

$ws.Range("D46").Value = "'89.06"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0960"
$ws.Range("E50").Value = "  -9.74%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.16%  "
